$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right above row 726, pushing the existing rows
# (old 726..775) down to (728..777). This mirrors the source diff, where
# every row from 728 downward simply carries the data that used to sit
# two rows higher, and the block ends up two rows longer (new last row 777).
$ws.Rows.Item(726).Resize(2).Insert()

# Populate the two newly-inserted rows with the new "Macroferia Regional
# de Talca" / "Platano" records.

# Row 726: Calidad "Pintón"
$ws.Range("A726").Value = 5
$ws.Range("B726").Value = "Macroferia Regional de Talca"
$ws.Range("C726").Value = "Maule"
$ws.Range("D726").Value = 44826
$ws.Range("E726").Value = 7
$ws.Range("F726").Value = "Fruta"
$ws.Range("G726").Value = 100108
$ws.Range("H726").Value = "Tropicales y subtropicales"
$ws.Range("I726").Value = 100108006
$ws.Range("J726").Value = "Plátano"
$ws.Range("K726").Value = "Sin especificar"
$ws.Range("L726").Value = "Pintón"
$ws.Range("M726").Value = 850
$ws.Range("N726").Value = 20000
$ws.Range("O726").Value = 20000
$ws.Range("P726").Value = 20000
$ws.Range("Q726").Value = "$/caja 20 kilos"
$ws.Range("R726").Value = "Ecuador"
$ws.Range("S726").Value = 1000
$ws.Range("T726").Value = 20

# Row 727: Calidad "Primera Pintón"
$ws.Range("A727").Value = 5
$ws.Range("B727").Value = "Macroferia Regional de Talca"
$ws.Range("C727").Value = "Maule"
$ws.Range("D727").Value = 44826
$ws.Range("E727").Value = 7
$ws.Range("F727").Value = "Fruta"
$ws.Range("G727").Value = 100108
$ws.Range("H727").Value = "Tropicales y subtropicales"
$ws.Range("I727").Value = 100108006
$ws.Range("J727").Value = "Plátano"
$ws.Range("K727").Value = "Sin especificar"
$ws.Range("L727").Value = "Primera Pintón"
$ws.Range("M727").Value = 560
$ws.Range("N727").Value = 21000
$ws.Range("O727").Value = 21000
$ws.Range("P727").Value = 21000
$ws.Range("Q727").Value = "$/caja 20 kilos"
$ws.Range("R727").Value = "Ecuador"
$ws.Range("S727").Value = 1050
$ws.Range("T727").Value = 20
